$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Restore the original Abstract/Authors text for the two rows whose
# shared-string references were hijacked by duplicated (whitespace-
# mangled) entries in the reverted commit.

$abstractRow19 = @"
Objective
id="Par1">To investigate the occurrence of olfactory and gustatory dysfunctions in patients with laboratory-confirmed COVID-19 infection.


Methods
id="Par2">Patients with laboratory-confirmed COVID-19 infection were recruited from 12 European hospitals.

 The following epidemiological and clinical outcomes have been studied: age, sex, ethnicity, comorbidities, and general and otolaryngological symptoms.

 Patients completed olfactory and gustatory questionnaires based on the smell and taste component of the National Health and Nutrition Examination Survey, and the short version of the Questionnaire of Olfactory Disorders-Negative Statements (sQOD-NS).


Results
id="Par3">A total of 417 mild-to-moderate COVID-19 patients completed the study (263 females).

 The most prevalent general symptoms consisted of cough, myalgia, and loss of appetite.

 Face pain and nasal obstruction were the most disease-related otolaryngological symptoms.

 85.6% and 88.0% of patients reported olfactory and gustatory dysfunctions, respectively.

 There was a significant association between both disorders (p &lt; 0.001).

 Olfactory dysfunction (OD) appeared before the other symptoms in 11.8% of cases.

 The sQO-NS scores were significantly lower in patients with anosmia compared with normosmic or hyposmic individuals (p = 0.001).

 Among the 18.2% of patients without nasal obstruction or rhinorrhea, 79.7% were hyposmic or anosmic.

 The early olfactory recovery rate was 44.0%.

 Females were significantly more affected by olfactory and gustatory dysfunctions than males (p = 0.001).


Conclusion
id="Par4">Olfactory and gustatory disorders are prevalent symptoms in European COVID-19 patients, who may not have nasal symptoms.

 The sudden anosmia or ageusia need to be recognized by the international scientific community as important symptoms of the COVID-19 infection.


Electronic supplementary material
The online version of this article (10.1007/s00405-020-05965-1) contains supplementary material, which is available to authorized users.



"@

$authorsRow19 = @"
[Jerome R.%Lechien%Jerome.Lechien@umons.ac.be%2,  Carlos M.%Chiesa-Estomba%NULL%2,  Daniele R.%De Siati%NULL%2,  Mihaela%Horoi%NULL%2,  Serge D.%Le Bon%NULL%2,  Alexandra%Rodriguez%NULL%2,  Didier%Dequanter%NULL%2,  Serge%Blecic%NULL%2,  Fahd%El Afia%NULL%2,  Lea%Distinguin%NULL%2,  Younes%Chekkoury-Idrissi%NULL%2,  Stéphane%Hans%NULL%2,  Irene Lopez%Delgado%NULL%2,  Christian%Calvo-Henriquez%NULL%2,  Philippe%Lavigne%NULL%2,  Chiara%Falanga%NULL%2,  Maria Rosaria%Barillari%NULL%2,  Giovanni%Cammaroto%NULL%2,  Mohamad%Khalife%NULL%2,  Pierre%Leich%NULL%2,  Christel%Souchay%NULL%2,  Camelia%Rossi%NULL%2,  Fabrice%Journe%NULL%2,  Julien%Hsieh%NULL%2,  Myriam%Edjlali%NULL%2,  Robert%Carlier%NULL%2,  Laurence%Ris%NULL%2,  Andrea%Lovato%NULL%2,  Cosimo%De Filippis%NULL%2,  Frederique%Coppee%NULL%2,  Nicolas%Fakhry%NULL%2,  Tareck%Ayad%NULL%2,  Sven%Saussez%NULL%2]
"@

$abstractRow26 = @"
Purpose
id="Par1">COVID-19 displays a variety of clinical manifestations; in pauci-symptomatic patients olfactory (OD) and gustatory dysfunctions (GD) may represent the first or only symptom.

 This topic is currently arousing great interest, and a growing number of papers are being published.

 Aim of this study is to investigate the timing of recovery from OD and GD in a real-life population hospitalized for COVID-19.
Methods
id="Par2">We followed up by a phone interview the first 100 patients discharged a month earlier from three Italian non-intensive care wards.


Results
id="Par3">All 100 patients were Caucasian, mean age was 65 years, 60% were males.

 Forty-two patients (mean age 63 years) experienced subjective chemosensory dysfunctions (29 OD and 41 GD): the male/female ratio was 2:1; 83% reported a complete or near complete recovery at follow-up.

 The recovery rate was not significantly different between males and females.

 The mean duration of OD and GD was 18 and 16 days, respectively.

 The mean recovery time from OD or GD resulted significantly longer for females than for males (26 vs 14 days, P = 0.009).

 Among the 42 symptomatic, the mean age of males was significantly higher than that of females (66 vs 57 years, P = 0.04), while the opposite was observed in the 58 asymptomatic patients (60 vs 73 years, P = 0.0018).


Conclusions
id="Par4">Recovery from OD or GD was rapid, occurring within 4 weeks in most patients.

 Chemosensory dysfunctions in women was less frequent, but longer lasting.

 The value of our study is its focus on a population of hospitalized patients significantly older than those previously described, and the additional data on gender differences.



"@

$authorsRow26 = @"
[Simone%Meini%simonemeini2@gmail.com%1,  Lorenzo Roberto%Suardi%NULL%2,  Lorenzo Roberto%Suardi%NULL%0,  Michele%Busoni%NULL%1,  Anna Teresa%Roberts%NULL%1,  Alberto%Fortini%NULL%1]
"@

$ws.Range("D19").Value2 = $abstractRow19
$ws.Range("E19").Value2 = $authorsRow19
$ws.Range("D26").Value2 = $abstractRow26
$ws.Range("E26").Value2 = $authorsRow26
